$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3119.6
$ws.Range("I64").Value = 2950
$ws.Range("J64").Value = 3232.6667
$ws.Range("K64").Value = 2950
$ws.Range("L64").Value = 3232.6667
$ws.Range("M64").Value = -2702
$ws.Range("N64").Value = -3728.6667
$ws.Range("H67").Value = 3119.6
$ws.Range("I67").Value = 2950
$ws.Range("J67").Value = 3232.6667
$ws.Range("K67").Value = 2950
$ws.Range("L67").Value = 3232.6667
$ws.Range("M67").Value = -2092
$ws.Range("N67").Value = -4948.6667
$ws.Range("H99").Value = 2641.25
$ws.Range("J99").Value = 6251
$ws.Range("L99").Value = 18753
$ws.Range("N99").Value = -21749
$ws.Range("H125").Value = 1237.5
$ws.Range("J125").Value = 1483.3334
$ws.Range("L125").Value = 13350.0006
$ws.Range("N125").Value = -18270.0006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2243.72
$ws.Range("I2").Value = 1509.85
$ws.Range("J2").Value = 5179.2
$ws.Range("K2").Value = 1509.85
$ws.Range("L2").Value = 5179.2
$ws.Range("M2").Value = -1396.85
$ws.Range("N2").Value = -5405.2
$ws.Range("H32").Value = 3050.077
$ws.Range("I32").Value = 1786.909
$ws.Range("K32").Value = 1786.909
$ws.Range("M32").Value = -1499.909
$ws.Range("H54").Value = 69025.4
$ws.Range("J54").Value = 69025.4
$ws.Range("L54").Value = 69025.4
$ws.Range("N54").Value = -70563.4
$ws.Range("H61").Value = 3895.625
$ws.Range("J61").Value = 4488.25
$ws.Range("L61").Value = 4488.25
$ws.Range("N61").Value = -4912.25
$ws.Range("H116").Value = 2243.72
$ws.Range("I116").Value = 1509.85
$ws.Range("J116").Value = 5179.2
$ws.Range("K116").Value = 1509.85
$ws.Range("L116").Value = 5179.2
$ws.Range("M116").Value = 784.1500000000001
$ws.Range("N116").Value = -9767.2
$ws.Range("H122").Value = 2932.5715
$ws.Range("I122").Value = 2793.3845
$ws.Range("J122").Value = 3158.75
$ws.Range("K122").Value = 8380.1535
$ws.Range("L122").Value = 9476.25
$ws.Range("M122").Value = -5930.1535
$ws.Range("N122").Value = -14376.25
$ws.Range("H132").Value = 7802.533
$ws.Range("I132").Value = 4162.6
$ws.Range("J132").Value = 15082.4
$ws.Range("K132").Value = 12487.8
$ws.Range("L132").Value = 45247.2
$ws.Range("M132").Value = -9957.800000000001
$ws.Range("N132").Value = -50307.2
$ws.Range("H136").Value = 3895.625
$ws.Range("J136").Value = 4488.25
$ws.Range("L136").Value = 13464.75
$ws.Range("N136").Value = -18564.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2243.72
$ws.Range("I3").Value = 1509.85
$ws.Range("J3").Value = 5179.2
$ws.Range("K3").Value = 1509.85
$ws.Range("L3").Value = 5179.2
$ws.Range("M3").Value = -1395.85
$ws.Range("N3").Value = -5407.2
$ws.Range("H134").Value = 3495.5789
$ws.Range("I134").Value = 3356.5557
$ws.Range("K134").Value = 10069.6671
$ws.Range("M134").Value = -7534.667099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1695.2858
$ws.Range("I31").Value = 1573.4
$ws.Range("K31").Value = 1573.4
$ws.Range("M31").Value = -1278.4
$ws.Range("H34").Value = 1695.2858
$ws.Range("I34").Value = 1573.4
$ws.Range("K34").Value = 1573.4
$ws.Range("M34").Value = -1371.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 474.66666
$ws.Range("H48").Value = 9999
$ws.Range("J48").Value = 9999
$ws.Range("L48").Value = 29997
$ws.Range("N48").Value = -30497
$ws.Range("H62").Value = 6999
$ws.Range("J62").Value = 6999
$ws.Range("L62").Value = 20997
$ws.Range("N62").Value = -22369
$ws.Range("H65").Value = 6999
$ws.Range("J65").Value = 6999
$ws.Range("L65").Value = 62991
$ws.Range("N65").Value = -69855
$ws.Range("H123").Value = 7497
$ws.Range("I123").Value = 7497
$ws.Range("K123").Value = 22491
$ws.Range("M123").Value = -20041
$ws.Range("H124").Value = 2299.5
$ws.Range("I124").Value = 2299.5
$ws.Range("K124").Value = 6898.5
$ws.Range("M124").Value = -1988.5
$ws.Range("H135").Value = 474.66666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2425
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 2200
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 2200
$ws.Range("M80").Value = -1502
$ws.Range("N80").Value = -4196
$ws.Range("H83").Value = 2425
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 2200
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 11000
$ws.Range("M83").Value = -7508
$ws.Range("N83").Value = -20984
$ws.Range("H126").Value = 7649.2144
$ws.Range("I126").Value = 6614.5713
$ws.Range("K126").Value = 19843.7139
$ws.Range("M126").Value = -17373.7139
$ws.Range("H132").Value = 5216.95
$ws.Range("I132").Value = 5502.278
$ws.Range("J132").Value = 2649
$ws.Range("K132").Value = 16506.834
$ws.Range("L132").Value = 7947
$ws.Range("M132").Value = -13976.834
$ws.Range("N132").Value = -13007
$ws.Range("H140").Value = 72571.336
$ws.Range("J140").Value = 72571.336
$ws.Range("L140").Value = 72571.336
$ws.Range("N140").Value = -82931.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 500015000
$ws.Range("I18").Value = 500015000
$ws.Range("K18").Value = 500015000
$ws.Range("M18").Value = -500014828
$ws.Range("H136").Value = 2819.4443
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 2910.7144
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 8732.143199999999
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -13832.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H122").Value = 5998.8237
$ws.Range("I122").Value = 6453.8
$ws.Range("J122").Value = 5348.857
$ws.Range("K122").Value = 19361.4
$ws.Range("L122").Value = 16046.571
$ws.Range("M122").Value = -16911.4
$ws.Range("N122").Value = -20946.571
